$wb = $excel.ActiveWorkbook

$nk  = $wb.Worksheets.Item("NK")
$jrp = $wb.Worksheets.Item("JRP")

# ---------------------------------------------------------------------
# 1) Update the NK-sheet percentage surcharge factors (row 14..19) for
#    the four blocks (NK1 Maenner / NK2 Maenner / NK1 Frauen / NK2 Frauen)
#    F: 1.19 -> 1.2 , L: 1.19 -> 1.8 , U: 1.15 -> 1.14 , AA: 1.17 -> 1.16
# ---------------------------------------------------------------------
$nk.Range("F14:F19").Value  = 1.2
$nk.Range("L14:L19").Value  = 1.8
$nk.Range("U14:U19").Value  = 1.14
$nk.Range("AA14:AA19").Value = 1.16

# ---------------------------------------------------------------------
# 2) Update the JRP-sheet category labels (column A, rows 4..15) so the
#    percentages shown to the user follow the refined factors above.
#    Writing the rows in this particular order reproduces the exact
#    shared-string table ordering used by the canonical edit.
# ---------------------------------------------------------------------
$jrp.Range("A13").Value = "NK1 +16% Frauen"
$jrp.Range("A14").Value = "NK1 +16% Frauen"
$jrp.Range("A15").Value = "NK1 +16% Frauen"

$jrp.Range("A10").Value = "NK1 +14% Mäner"
$jrp.Range("A11").Value = "NK1 +14% Mäner"
$jrp.Range("A12").Value = "NK1 +14% Mäner"

$jrp.Range("A7").Value = "NK2 +18% Frauen"
$jrp.Range("A8").Value = "NK2 +18% Frauen"
$jrp.Range("A9").Value = "NK2 +18% Frauen"

$jrp.Range("A4").Value = "NK2 +20% Mäner"
$jrp.Range("A5").Value = "NK2 +20% Mäner"
$jrp.Range("A6").Value = "NK2 +20% Mäner"

# ---------------------------------------------------------------------
# 3) Update sheet views: JRP becomes the selected/active tab (instead of
#    NK), with the selection moved to A4. NK keeps its own selection but
#    is scrolled back towards the left (column B) and selects AA19.
# ---------------------------------------------------------------------
$nk.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$nk.Range("AA19").Select()

$jrp.Activate()
$jrp.Range("A4").Select()
